# Apply "custom accuracy" rounding (2 decimal places) to row 5 data cells (B5:AH5),
# then remove the last data row (row 6), which reduces the sheet's used range
# from A1:AH6 down to A1:AH5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH")

foreach ($col in $cols) {
    $addr = $col + "5"
    $cur = $ws.Range($addr).Value2
    $rounded = $excel.WorksheetFunction.Round($cur, 2)
    $ws.Range($addr).Value = $rounded
}

# Delete row 6 entirely (shifts dimension down to A1:AH5)
$ws.Rows.Item(6).Delete()
